$d = $word.ActiveDocument

$find = $d.Content.Find
$find.Text = "100: _____ ("
$find.Replacement.Text = "100 knots: _____ ("
$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)
